# Scheduled runner update: refresh market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) for the Leve
# profit tables across all Disciple of the Hand job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 637.7778
$ws.Range("J17").Value = 659.88464
$ws.Range("L17").Value = 1979.65392
$ws.Range("N17").Value = -2315.65392

$ws.Range("H76").Value = 16675115
$ws.Range("I76").Value = 9379
$ws.Range("K76").Value = 9379
$ws.Range("M76").Value = -9064

$ws.Range("H79").Value = 16675115
$ws.Range("I79").Value = 9379
$ws.Range("K79").Value = 9379
$ws.Range("M79").Value = -8287

$ws.Range("H86").Value = 50266300
$ws.Range("I86").Value = 66668388
$ws.Range("J86").Value = 9261075
$ws.Range("K86").Value = 66668388
$ws.Range("L86").Value = 9261075
$ws.Range("M86").Value = -66667265
$ws.Range("N86").Value = -9263321

$ws.Range("H87").Value = 59997
$ws.Range("J87").Value = 59997
$ws.Range("L87").Value = 59997
$ws.Range("N87").Value = -62493

$ws.Range("H88").Value = 78762040
$ws.Range("J88").Value = 11181110
$ws.Range("L88").Value = 11181110
$ws.Range("N88").Value = -11181922

$ws.Range("H89").Value = 50266300
$ws.Range("I89").Value = 66668388
$ws.Range("J89").Value = 9261075
$ws.Range("K89").Value = 333341940
$ws.Range("L89").Value = 46305375
$ws.Range("M89").Value = -333336324
$ws.Range("N89").Value = -46316607

$ws.Range("H90").Value = 59997
$ws.Range("J90").Value = 59997
$ws.Range("L90").Value = 179991
$ws.Range("N90").Value = -192471

$ws.Range("H91").Value = 78762040
$ws.Range("J91").Value = 11181110
$ws.Range("L91").Value = 11181110
$ws.Range("N91").Value = -11183918

$ws.Range("H98").Value = 9569
$ws.Range("I98").Value = 9569
$ws.Range("K98").Value = 9569
$ws.Range("M98").Value = -8071

$ws.Range("H112").Value = 5652.59
$ws.Range("J112").Value = 6155.1143
$ws.Range("L112").Value = 18465.3429
$ws.Range("N112").Value = -20681.3429

$ws.Range("H122").Value = 9569
$ws.Range("I122").Value = 9569
$ws.Range("K122").Value = 28707
$ws.Range("M122").Value = -26257

$ws.Range("H137").Value = 3001.0952
$ws.Range("I137").Value = 3280.5557
$ws.Range("K137").Value = 9841.667099999999
$ws.Range("M137").Value = -7291.667099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2121543.5
$ws.Range("I32").Value = 2406674.5
$ws.Range("K32").Value = 2406674.5
$ws.Range("M32").Value = -2406387.5

$ws.Range("H88").Value = 125002000
$ws.Range("I88").Value = 4000
$ws.Range("J88").Value = 250000000
$ws.Range("K88").Value = 4000
$ws.Range("L88").Value = 250000000
$ws.Range("M88").Value = -3594
$ws.Range("N88").Value = -250000812

$ws.Range("H91").Value = 125002000
$ws.Range("I91").Value = 4000
$ws.Range("J91").Value = 250000000
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 250000000
$ws.Range("M91").Value = -2596
$ws.Range("N91").Value = -250002808

$ws.Range("H122").Value = 2875.2327
$ws.Range("I122").Value = 2314.5312
$ws.Range("J122").Value = 4506.364
$ws.Range("K122").Value = 6943.5936
$ws.Range("L122").Value = 13519.092
$ws.Range("M122").Value = -4493.5936
$ws.Range("N122").Value = -18419.092

$ws.Range("H132").Value = 4423.125
$ws.Range("I132").Value = 1971.4814
$ws.Range("J132").Value = 9515
$ws.Range("K132").Value = 5914.4442
$ws.Range("L132").Value = 28545
$ws.Range("M132").Value = -3384.4442
$ws.Range("N132").Value = -33605

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 355.7857
$ws.Range("I94").Value = 218.26315
$ws.Range("J94").Value = 646.1111
$ws.Range("K94").Value = 218.26315
$ws.Range("L94").Value = 646.1111
$ws.Range("M94").Value = 232.73685
$ws.Range("N94").Value = -1548.1111

$ws.Range("H99").Value = 2163.5
$ws.Range("I99").Value = 871.75
$ws.Range("J99").Value = 4747
$ws.Range("K99").Value = 871.75
$ws.Range("L99").Value = 4747
$ws.Range("M99").Value = 626.25
$ws.Range("N99").Value = -7743

$ws.Range("H134").Value = 8061.607
$ws.Range("I134").Value = 4011
$ws.Range("K134").Value = 12033
$ws.Range("M134").Value = -9498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 499.66666
$ws.Range("I22").Value = 499.66666
$ws.Range("K22").Value = 499.66666
$ws.Range("M22").Value = -149.66666

$ws.Range("H31").Value = 10420.84
$ws.Range("I31").Value = 4627
$ws.Range("J31").Value = 15769
$ws.Range("K31").Value = 4627
$ws.Range("L31").Value = 15769
$ws.Range("M31").Value = -4332
$ws.Range("N31").Value = -16359

$ws.Range("H34").Value = 10420.84
$ws.Range("I34").Value = 4627
$ws.Range("J34").Value = 15769
$ws.Range("K34").Value = 4627
$ws.Range("L34").Value = 15769
$ws.Range("M34").Value = -4425
$ws.Range("N34").Value = -16173

$ws.Range("H132").Value = 7428.5
$ws.Range("I132").Value = 4491.273
$ws.Range("J132").Value = 9913.846
$ws.Range("K132").Value = 13473.819
$ws.Range("L132").Value = 29741.538
$ws.Range("M132").Value = -10943.819
$ws.Range("N132").Value = -34801.538

$ws.Range("H134").Value = 11196.634
$ws.Range("I134").Value = 16785.555
$ws.Range("J134").Value = 8801.380999999999
$ws.Range("K134").Value = 50356.665
$ws.Range("L134").Value = 26404.143
$ws.Range("M134").Value = -47821.665
$ws.Range("N134").Value = -31474.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 95
$ws.Range("I19").Value = 95
$ws.Range("K19").Value = 285
$ws.Range("M19").Value = -111

$ws.Range("H68").Value = 7167.143
$ws.Range("J68").Value = 11305.25
$ws.Range("L68").Value = 33915.75
$ws.Range("N68").Value = -35537.75

$ws.Range("H71").Value = 7167.143
$ws.Range("J71").Value = 11305.25
$ws.Range("L71").Value = 101747.25
$ws.Range("N71").Value = -109859.25

$ws.Range("H107").Value = 846.5
$ws.Range("I107").Value = 499
$ws.Range("K107").Value = 1497
$ws.Range("M107").Value = 423

$ws.Range("H113").Value = 3068.037
$ws.Range("I113").Value = 667.2222
$ws.Range("J113").Value = 4268.4443
$ws.Range("K113").Value = 2001.6666
$ws.Range("L113").Value = 12805.3329
$ws.Range("M113").Value = 168.3334
$ws.Range("N113").Value = -17145.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 33335434
$ws.Range("J126").Value = 2332.6667
$ws.Range("L126").Value = 6998.000100000001
$ws.Range("N126").Value = -11938.0001

$ws.Range("H132").Value = 4959.9473
$ws.Range("I132").Value = 1296.4546
$ws.Range("J132").Value = 9997.25
$ws.Range("K132").Value = 3889.3638
$ws.Range("L132").Value = 29991.75
$ws.Range("M132").Value = -1359.3638
$ws.Range("N132").Value = -35051.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4318.364
$ws.Range("I7").Value = 3071.4285
$ws.Range("K7").Value = 3071.4285
$ws.Range("M7").Value = -2959.4285

$ws.Range("H16").Value = 1279.619
$ws.Range("I16").Value = 1279.619
$ws.Range("K16").Value = 1279.619
$ws.Range("M16").Value = -1109.619

$ws.Range("H61").Value = 3459.641
$ws.Range("I61").Value = 2476.1365
$ws.Range("J61").Value = 4732.4116
$ws.Range("K61").Value = 2476.1365
$ws.Range("L61").Value = 4732.4116
$ws.Range("M61").Value = -2274.1365
$ws.Range("N61").Value = -5136.4116

$ws.Range("H68").Value = 166670190
$ws.Range("J68").Value = 5132.6665
$ws.Range("L68").Value = 5132.6665
$ws.Range("N68").Value = -6630.6665

$ws.Range("H71").Value = 166670190
$ws.Range("J71").Value = 5132.6665
$ws.Range("L71").Value = 25663.3325
$ws.Range("N71").Value = -33151.3325

$ws.Range("H113").Value = 3459.641
$ws.Range("I113").Value = 2476.1365
$ws.Range("J113").Value = 4732.4116
$ws.Range("K113").Value = 2476.1365
$ws.Range("L113").Value = 4732.4116
$ws.Range("M113").Value = -306.1365000000001
$ws.Range("N113").Value = -9072.411599999999

$ws.Range("H122").Value = 3370.6216
$ws.Range("I122").Value = 2784.0435
$ws.Range("J122").Value = 4334.2856
$ws.Range("K122").Value = 8352.130500000001
$ws.Range("L122").Value = 13002.8568
$ws.Range("M122").Value = -5902.130500000001
$ws.Range("N122").Value = -17902.8568

$ws.Range("H126").Value = 4318.364
$ws.Range("I126").Value = 3071.4285
$ws.Range("K126").Value = 9214.2855
$ws.Range("M126").Value = -6744.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8800
$ws.Range("I62").Value = 8083.6665
$ws.Range("J62").Value = 9874.5
$ws.Range("K62").Value = 8083.6665
$ws.Range("L62").Value = 9874.5
$ws.Range("M62").Value = -7459.6665
$ws.Range("N62").Value = -11122.5

$ws.Range("H65").Value = 8800
$ws.Range("I65").Value = 8083.6665
$ws.Range("J65").Value = 9874.5
$ws.Range("K65").Value = 40418.3325
$ws.Range("L65").Value = 49372.5
$ws.Range("M65").Value = -37298.3325
$ws.Range("N65").Value = -55612.5
